# Refined metadata to be additional tab
$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1) Update the "time_taken" timestamps on the data sheet (column F, rows 2-8) ---
$dataSheet.Range("F2").Value = "2021-10-05 14:20:07.693224"
$dataSheet.Range("F3").Value = "2021-10-05 14:20:07.693232"
$dataSheet.Range("F4").Value = "2021-10-05 14:20:07.693235"
$dataSheet.Range("F5").Value = "2021-10-05 14:20:07.693238"
$dataSheet.Range("F6").Value = "2021-10-05 14:20:07.693240"
$dataSheet.Range("F7").Value = "2021-10-05 14:20:07.693243"
$dataSheet.Range("F8").Value = "2021-10-05 14:20:07.693246"

# --- 2) Add a new "metadata" worksheet right after "data" ---
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row (row 1) in B1:G1, formatted like the "data" sheet's header row (bold/border/center - style index 1)
$headers = @("data_name", "data_id", "data_version", "data_version_created", "panel_query_time", "panel_get_request")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row (row 2)
$ws.Cells.Item(2, 1).Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Cells.Item(2, 2).Value = "Epidermodysplasia verruciformis"
$ws.Cells.Item(2, 3).Value = 562

# "1.4" needs to remain text, not be coerced into a number - force text then strip the
# resulting style back to the default (unstyled) by pasting formats from an unstyled cell.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "1.4"
$dataSheet.Range("B2").Copy()
$ws.Range("D2").PasteSpecial(-4122)

$ws.Cells.Item(2, 5).Value = "2020-10-15T19:09:35.918741Z"
$ws.Cells.Item(2, 6).Value = "2021-10-05 14:20:07.689508"
$ws.Cells.Item(2, 7).Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/562/?format=json"

Write-Output "done"
